# Apply scheduled-runner profit recalculation updates across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 272.92856
$ws.Range("I33").Value = 178.53847
$ws.Range("K33").Value = 178.53847
$ws.Range("M33").Value = 50.46153000000001
$ws.Range("H40").Value = 5676.773
$ws.Range("I40").Value = 1466.6666
$ws.Range("K40").Value = 1466.6666
$ws.Range("M40").Value = -1291.6666
$ws.Range("H64").Value = 4800
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("M64").Value = -4752
$ws.Range("H67").Value = 4800
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("M67").Value = -4142
$ws.Range("H98").Value = 738.63635
$ws.Range("I98").Value = 738.63635
$ws.Range("K98").Value = 738.63635
$ws.Range("M98").Value = 759.36365
$ws.Range("H111").Value = 625.625
$ws.Range("I111").Value = 502.66666
$ws.Range("K111").Value = 1507.99998
$ws.Range("M111").Value = 1559.00002
$ws.Range("H113").Value = 6022.5
$ws.Range("I113").Value = 4780.7
$ws.Range("K113").Value = 4780.7
$ws.Range("M113").Value = -1526.7
$ws.Range("H122").Value = 738.63635
$ws.Range("I122").Value = 738.63635
$ws.Range("K122").Value = 2215.90905
$ws.Range("M122").Value = 234.0909499999998
$ws.Range("H138").Value = 1735.2307
$ws.Range("I138").Value = 944.75
$ws.Range("K138").Value = 2834.25
$ws.Range("M138").Value = 2305.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 697.5
$ws.Range("I2").Value = 697.5
$ws.Range("K2").Value = 697.5
$ws.Range("M2").Value = -584.5
$ws.Range("H74").Value = 3357.5715
$ws.Range("I74").Value = 3198.1667
$ws.Range("K74").Value = 3198.1667
$ws.Range("M74").Value = -2324.1667
$ws.Range("H77").Value = 3357.5715
$ws.Range("I77").Value = 3198.1667
$ws.Range("K77").Value = 15990.8335
$ws.Range("M77").Value = -11622.8335
$ws.Range("H102").Value = 3000
$ws.Range("J102").Value = 3000
$ws.Range("L102").Value = 3000
$ws.Range("N102").Value = -6244
$ws.Range("H116").Value = 697.5
$ws.Range("I116").Value = 697.5
$ws.Range("K116").Value = 697.5
$ws.Range("M116").Value = 1596.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 697.5
$ws.Range("I3").Value = 697.5
$ws.Range("K3").Value = 697.5
$ws.Range("M3").Value = -583.5
$ws.Range("H10").Value = 1687.5
$ws.Range("I10").Value = 1375
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 1375
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = -1235
$ws.Range("N10").Value = -2280
$ws.Range("H18").Value = 26000
$ws.Range("J18").Value = 26000
$ws.Range("L18").Value = 26000
$ws.Range("N18").Value = -27058
$ws.Range("H25").Value = 1271.75
$ws.Range("I25").Value = 1271.75
$ws.Range("K25").Value = 1271.75
$ws.Range("M25").Value = -1036.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1917.75
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 2223.6667
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 2223.6667
$ws.Range("M2").Value = -887
$ws.Range("N2").Value = -2449.6667
$ws.Range("H3").Value = 2747
$ws.Range("J3").Value = 2747
$ws.Range("L3").Value = 2747
$ws.Range("N3").Value = -2973
$ws.Range("H132").Value = 4250
$ws.Range("I132").Value = 4250
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12750
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10220
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1650
$ws.Range("J46").Value = 2600
$ws.Range("L46").Value = 7800
$ws.Range("N46").Value = -7982
$ws.Range("H139").Value = 1833
$ws.Range("I139").Value = 1500
$ws.Range("K139").Value = 4500
$ws.Range("M139").Value = 640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 2000
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -2224
$ws.Range("H33").Value = 14666
$ws.Range("J33").Value = 14666
$ws.Range("L33").Value = 14666
$ws.Range("N33").Value = -15170
$ws.Range("H57").Value = 20000
$ws.Range("J57").Value = 20000
$ws.Range("L57").Value = 20000
$ws.Range("N57").Value = -21640
$ws.Range("H80").Value = 2771.1428
$ws.Range("I80").Value = 2279.6
$ws.Range("K80").Value = 2279.6
$ws.Range("M80").Value = -1281.6
$ws.Range("H83").Value = 2771.1428
$ws.Range("I83").Value = 2279.6
$ws.Range("K83").Value = 11398
$ws.Range("M83").Value = -6406
$ws.Range("H97").Value = 603.1667
$ws.Range("I97").Value = 603.1667
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 603.1667
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -107.1667
$ws.Range("N97").ClearContents()
$ws.Range("H107").Value = 272.4
$ws.Range("I107").Value = 272.4
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 272.4
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1647.6
$ws.Range("N107").ClearContents()
$ws.Range("H122").Value = 1046
$ws.Range("I122").Value = 1051.4
$ws.Range("K122").Value = 3154.2
$ws.Range("M122").Value = -704.2000000000003
$ws.Range("H126").Value = 12435
$ws.Range("J126").Value = 15499.667
$ws.Range("L126").Value = 46499.001
$ws.Range("N126").Value = -51439.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19730.1
$ws.Range("J7").Value = 26625
$ws.Range("L7").Value = 26625
$ws.Range("N7").Value = -26849
$ws.Range("H40").Value = 4217.2104
$ws.Range("I40").Value = 3895.9443
$ws.Range("K40").Value = 3895.9443
$ws.Range("M40").Value = -3759.9443
$ws.Range("H61").Value = 2282.8333
$ws.Range("I61").Value = 2282.8333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2282.8333
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2080.8333
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 2282.8333
$ws.Range("I113").Value = 2282.8333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2282.8333
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -112.8332999999998
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 19730.1
$ws.Range("J126").Value = 26625
$ws.Range("L126").Value = 79875
$ws.Range("N126").Value = -84815

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 104000
$ws.Range("J2").Value = 104000
$ws.Range("L2").Value = 104000
$ws.Range("N2").Value = -104224
$ws.Range("H11").Value = 41725000
$ws.Range("J11").Value = 87500
$ws.Range("L11").Value = 87500
$ws.Range("N11").Value = -87784
$ws.Range("H122").Value = 2499.2
$ws.Range("I122").Value = 2499.2
$ws.Range("K122").Value = 7497.599999999999
$ws.Range("M122").Value = -5047.599999999999
$ws.Range("H132").Value = 2516
$ws.Range("I132").Value = 2378
$ws.Range("K132").Value = 7134
$ws.Range("M132").Value = -4604
$ws.Range("H136").Value = 4047.3845
$ws.Range("I136").Value = 4578.8184
$ws.Range("J136").Value = 1124.5
$ws.Range("K136").Value = 13736.4552
$ws.Range("L136").Value = 3373.5
$ws.Range("M136").Value = -11186.4552
$ws.Range("N136").Value = -8473.5
